$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title / author ---
Replace-Text "Rhythm of Faith" "Mathematics - The Language of the Universe"
Replace-Text "Ali Hassan" "Henry Richardson"

# --- Email/handle paragraph: collapse 5 runs into a single run "at" ---
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3e = $d.Range($r3.Start, $r3.End - 1)
$r3e.Text = "at"

# --- Body paragraph (the long one with many sentence runs) ---
Replace-Text "In the heart of existence, a tapestry of faith unfurls, its radiant threads weaving a symphony of devotion" "Step into the fascinating world of mathematics, a language that unveils the secrets of the universe"
Replace-Text " From the sacred chants of ancient temples to the fervent prayers whispered in quiet corners, humanity's pursuit of the divine has given rise to a kaleidoscope of spiritual practices and traditions" " Since time immemorial, humans have used this versatile tool to decode patterns, make predictions, and solve complex problems"
Replace-Text " Across cultures and epochs, the rhythm of faith has pulsed through generations, shaping individual lives and collective destinies" " Beyond its practical applications, mathematics offers a profound understanding of the world we live in, revealing the intricate connections between seemingly disparate concepts"
Replace-Text " It is a force that transcends boundaries, uniting souls in a quest for meaning and purpose" " From counting pebbles to unraveling the enigma of quantum mechanics, mathematics has been an indispensable companion to humanity's journey of knowledge"
Replace-Text "Faith, like a beacon in the dark, guides us through the uncertainties of life" "Dive into the beauty of mathematical structures, where numbers, shapes, and equations intertwine in a harmonious dance"
Replace-Text " It provides solace in times of despair and ignites hope amidst adversity" " Discover the power of mathematical thinking, a skill that enables us to analyze, reason, and make informed decisions"
Replace-Text " It instills within us a sense of belonging, connecting us to a community of believers who share our aspirations and values" " Explore the rich tapestry of mathematical ideas, from the ancient wisdom of Euclid's geometry to the elegant simplicity of calculus"
Replace-Text " Whether it is the unwavering belief in a higher power, the reverence for sacred texts, or the observance of rituals and traditions, faith becomes a compass that navigates the intricate labyrinth of human existence" " In this realm of numbers and logic, there's a universe waiting to be explored - a universe where patterns emerge from chaos and order reigns supreme"
Replace-Text "Yet, faith is not a monolithic entity" "Delve into the depths of mathematical applications, where numbers and equations guide us through the intricacies of science, technology, engineering, and even art"
Replace-Text " It manifests in myriad forms, as diverse as the human spirit itself" " From calculating rocket trajectories to designing bridges, from predicting weather patterns to creating computer algorithms, mathematics is an indispensable tool in shaping our modern world"
Replace-Text " From the grand cathedrals of Christendom to the humble prayer mats of the faithful, from the meditative silence of Zen monasteries to the exuberant drumming of African tribal ceremonies, faith takes on countless expressions, each reflecting the unique cultural and historical context in which it flourishes" " It is the language of innovation, the engine of progress, and the key to unlocking the secrets of the future"

# --- Summary heading paragraph stays the same ("Summary") ---

# --- Summary body paragraph ---
Replace-Text "In the realm of human experience, faith emerges as a potent force, shaping lives and leaving an enduring imprint on the course of history" "Mathematics, the language of the universe, offers a profound understanding of the world around us"
Replace-Text " It transcends boundaries, uniting individuals in a shared quest for meaning and purpose" " It reveals intricate patterns, enables logical reasoning, and underpins numerous applications across various disciplines"

# Merge the final three sentence-runs (+ 2 period runs) of the Summary body paragraph into one run.
$summaryPara = $d.Paragraphs($d.Paragraphs.Count)
$sr = $summaryPara.Range
$needle = " Faith manifests in a kaleidoscope of expressions, ranging from the grand cathedrals of Christendom to the humble prayer mats of the faithful"
$startPos = $sr.Text.IndexOf($needle)
$mergeStart = $sr.Start + $startPos
$mergeEnd = $sr.End - 1  # exclude the trailing "." run and the paragraph mark
$mergeRange = $d.Range($mergeStart, $mergeEnd)
$mergeRange.Text = " From unraveling the mysteries of quantum mechanics to driving technological innovations, mathematics continues to be an indispensable tool for humanity's quest for knowledge and progress"

# --- Add a trailing empty paragraph at the very end of the document ---
$d.Content.Find.Execute("progress.", $true, $false, $false, $false, $false, $true, 1, $false, "progress.^p", 2) | Out-Null

Write-Output "done"
